$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix capitalization typo in the ExplorerType column: "FireFox" -> "Firefox"
$ws.Range("E2").Value = "Firefox"
